$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (model_four) updates
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 0.4558551907539368
$ws.Range("E5").Value = 0.8417368531227112
$ws.Range("F5").Value = 0.4799478054046631
$ws.Range("G5").Value = 0.8231666684150696
$ws.Range("H5").Value = 0.8394736842105263
$ws.Range("I5").Value = 0.8394736842105263
$ws.Range("J5").Value = 0.8394736842105263
$ws.Range("K5").Value = 0.8394736842105263
$ws.Range("L5").Value = "7:27"

# Row 8 (model_seven) updates
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 0.3552636504173279
$ws.Range("E8").Value = 0.8744210600852966
$ws.Range("F8").Value = 0.397102952003479
$ws.Range("G8").Value = 0.8615000247955322
$ws.Range("H8").Value = 0.8597368421052631
$ws.Range("I8").Value = 0.8597368421052631
$ws.Range("J8").Value = 0.8597368421052631
$ws.Range("K8").Value = 0.8597368421052631
$ws.Range("L8").Value = "14:38"
